# Apply the "update sample and main" edit:
#  - Fill in attendance rows 5-32 on Sheet1 with uuid/name/date data
#    (mirrors the pattern already present in rows 2-4), leaving the
#    timeIn/timeOut columns blank on the four rows that were blank
#    in the target (8, 16, 22, 29).
#  - Move the active sheet/selection from RD (sheet2) to Sheet1,
#    selecting D8:E8 there.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

$uuid = "ID1"
$name = "Alice Smith"

# Rows that keep timeIn (D) / timeOut (E) blank, matching the source diff.
$blankDE = @(8, 16, 22, 29)

for ($i = 0; $i -lt 28; $i++) {
    $row = 5 + $i
    $base = 45508 + $i

    $ws1.Cells.Item($row, 1).Value = $uuid          # A: uuid
    $ws1.Cells.Item($row, 2).Value = $name           # B: name
    $ws1.Cells.Item($row, 3).Value = $base + 0.375   # C: workingTime

    if ($blankDE -notcontains $row) {
        $ws1.Cells.Item($row, 4).Value = $base + 0.375          # D: timeIn
        $ws1.Cells.Item($row, 5).Value = $base + 0.7083333333   # E: timeOut
    }
}

# Switch the active tab from RD to Sheet1 and select D8:E8 there.
$ws1.Activate()
$ws1.Range("D8:E8").Select()
